$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correção no nome da coluna Regiao: era "TB Regiao", agora é "Regiao"
$ws.Range("C1").Value = "Regiao"

# Selecione a célula C1 (refletindo a seleção ativa após a edição)
$ws.Range("C1").Select()
